$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.553202
$ws.Range("H3").Value = 0.553202
$ws.Range("H4").Value = 0.553202
$ws.Range("H5").Value = 0.553202
$ws.Range("H6").Value = 0.553202

$ws.Range("M2").Value = 1.802415666666667
$ws.Range("N2").Value = 5.407247
$ws.Range("O2").Value = 0.1831574081839677
$ws.Range("P2").Value = 0.1961662442954491
$ws.Range("Q2").Value = 0.3323666505437777
$ws.Range("R2").Value = 2.991299854894
$ws.Range("S2").Value = 0.1831574081839677
$ws.Range("T2").Value = 0.1961662442954491

$ws.Range("O3").Value = 0.3425591289923409
$ws.Range("P3").Value = 0.3668895429883566
$ws.Range("Q3").Value = 0.6216250352375555
$ws.Range("R3").Value = 5.594625317137999
$ws.Range("S3").Value = 0.3425591289923409
$ws.Range("T3").Value = 0.3668895429883566

$ws.Range("M4").Value = 0.9204736666666666
$ws.Range("N4").Value = 2.761421
$ws.Range("O4").Value = 0.09353645455160088
$ws.Range("P4").Value = 0.1001799227016231
$ws.Range("Q4").Value = 0.1697359577824444
$ws.Range("R4").Value = 1.527623620042
$ws.Range("S4").Value = 0.09353645455160088
$ws.Range("T4").Value = 0.1001799227016231

$ws.Range("M5").Value = 1.957789
$ws.Range("N5").Value = 3.915578
$ws.Range("O5").Value = 0.1989460953112084
$ws.Range("P5").Value = 0.142050886616773
$ws.Range("Q5").Value = 0.3610175967926667
$ws.Range("R5").Value = 2.166105580756
$ws.Range("S5").Value = 0.1989460953112084
$ws.Range("T5").Value = 0.142050886616773

$ws.Range("M6").Value = 1.789066666666667
$ws.Range("N6").Value = 5.3672
$ws.Range("O6").Value = 0.1818009129608822
$ws.Range("P6").Value = 0.1947134033977982
$ws.Range("Q6").Value = 0.3299050860444444
$ws.Range("R6").Value = 2.969145774399999
$ws.Range("S6").Value = 0.1818009129608822
$ws.Range("T6").Value = 0.1947134033977982
